$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A155").Value = "Mümmel"
$ws.Range("B155").Value = "Phia"
$ws.Range("C155").Value = "Zoë"
$ws.Range("D155").Value = "Lieblingswort Zaza"
$ws.Range("E155").Value = 45458
$ws.Range("F155").Value = 5
$ws.Range("G155").Value = 5

$ws.Range("A156").Value = "Marian"
$ws.Range("B156").Value = "Emilio"
$ws.Range("C156").Value = "André"
$ws.Range("D156").Value = "Normale Kartoffeln auf die Eins"
$ws.Range("E156").Value = 45458
$ws.Range("F156").Value = 3
$ws.Range("G156").Value = 5

$ws.Range("A157").Value = "Verena "
$ws.Range("B157").Value = "Merlin"
$ws.Range("C157").Value = "Luca"
$ws.Range("D157").Value = "Pain in the Ass (Ace)"
$ws.Range("E157").Value = 45458
$ws.Range("F157").Value = 4
$ws.Range("G157").Value = 5

$ws.Range("A158").Value = "Marwin"
$ws.Range("B158").Value = "Friedrich"
$ws.Range("C158").Value = "Marie"
$ws.Range("D158").Value = "MFM"
$ws.Range("E158").Value = 45458
$ws.Range("F158").Value = 2
$ws.Range("G158").Value = 5

$ws.Range("A159").Value = "Chris"
$ws.Range("B159").Value = "Lorenz"
$ws.Range("C159").Value = "Valdemar"
$ws.Range("D159").Value = "Pferde auf die Eins"
$ws.Range("E159").Value = 45458
$ws.Range("F159").Value = 1
$ws.Range("G159").Value = 5

$ws.Range("E154").Copy()
$ws.Range("E155:E159").PasteSpecial(-4122)

[void]$ws.Range("H159").Select()
